$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("B2").Value = '2024-05-18'
$ws1.Range("C2").Value = '合肥·WA二次元饭局（取消）'
$ws1.Range("D2").Value = '临泉路胜利路交叉路（中环国际大厦对面） 太太满庭芳(胜利路店)'
$ws1.Range("E2").Value = '2024.05.18 14:50-05.18 20:00'
$ws1.Range("F2").Value = 67
$ws1.Range("G2").Value = '不可售'
$ws1.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=83978'
$ws1.Range("I2").Value = '//i1.hdslb.com/bfs/openplatform/202404/wK9Yq9Ta1712657384067.jpeg'

$ws1.Range("B3").Value = '2024-05-18'
$ws1.Range("C3").Value = '合肥·梦时空SPO1动漫展（取消）'
$ws1.Range("D3").Value = '阜阳路16号 银瑞林国际大酒店'
$ws1.Range("E3").Value = '2024.05.18 10:00-05.18 17:00'
$ws1.Range("F3").Value = 127
$ws1.Range("G3").Value = '不可售'
$ws1.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
$ws1.Range("I3").Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'

$ws1.Range("B4").Value = '2024-05-18'
$ws1.Range("C4").Value = '合肥·首届偶活企划——偶像计划-闪耀舞台'
$ws1.Range("D4").Value = '阜阳路16号 银瑞林国际大酒店'
$ws1.Range("E4").Value = '2024.05.18 09:00-05.18 17:00'
$ws1.Range("F4").Value = 60
$ws1.Range("G4").Value = 58
$ws1.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=83891'
$ws1.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202404/lfqv8l9Q1712453982625.jpeg'

$ws1.Range("B5").Value = '2024-06-01'
$ws1.Range("C5").Value = '合肥·运动番only·群青日和'
$ws1.Range("D5").Value = '金寨路287号 合肥明星运动公园'
$ws1.Range("E5").Value = '2024.06.01 09:30-06.01 17:30'
$ws1.Range("F5").Value = 504
$ws1.Range("G5").Value = 70
$ws1.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=83058'
$ws1.Range("I5").Value = '//i2.hdslb.com/bfs/openplatform/202404/Jzeq47lD1714026878824.jpeg'

$ws1.Range("B6").Value = '2024-06-08'
$ws1.Range("C6").Value = '合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~'
$ws1.Range("D6").Value = '锦绣大道3899号 合肥滨湖会展中心'
$ws1.Range("E6").Value = '2024.06.08 09:30-06.09 17:00'
$ws1.Range("F6").Value = 6740
$ws1.Range("G6").Value = 65
$ws1.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=83518'
$ws1.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202403/1Sqp42gM1711691520194.jpeg'

$ws1.Range("B7").Value = '2024-06-09'
$ws1.Range("C7").Value = '合肥·第二届华盟动漫次元嘉年华'
$ws1.Range("D7").Value = '常青街道十五里河村合柴1972院内 合肥当代美术馆'
$ws1.Range("E7").Value = '2024.06.09 10:00-06.10 17:00'
$ws1.Range("F7").Value = 186
$ws1.Range("G7").Value = 58
$ws1.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=84081'
$ws1.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202404/O5LyHE7j1712732240786.jpeg'

$ws1.Range("B8").Value = '2024-06-09'
$ws1.Range("C8").Value = '合肥·第六届环形宇宙动漫游戏嘉年华内场票·赵成晨'
$ws1.Range("D8").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws1.Range("E8").Value = '2024.06.09 09:30-06.09 17:00'
$ws1.Range("F8").Value = 148
$ws1.Range("G8").Value = 238
$ws1.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=84863'
$ws1.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202404/I5S4Ih2M1714031127805.jpeg'

$ws1.Range("B9").Value = '2024-06-22'
$ws1.Range("C9").Value = '合肥·Look Look动漫嘉年华'
$ws1.Range("D9").Value = '新站区东方大道288号 少荃体育中心'
$ws1.Range("E9").Value = '2024.06.22 10:00-06.22 17:30'
$ws1.Range("F9").Value = 1035
$ws1.Range("G9").Value = 58
$ws1.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=82311'
$ws1.Range("I9").Value = '//i2.hdslb.com/bfs/openplatform/202404/RFYwkzvt1713951750482.jpeg'

$ws1.Range("B10").Value = '2024-06-22'
$ws1.Range("C10").Value = '合肥·城市动漫节'
$ws1.Range("D10").Value = '包河经济开发区大连路与园博大道交口骆岗中央公园园博小镇一期S6区1号楼 大机库演艺中心'
$ws1.Range("E10").Value = '2024.06.22 10:00-06.23 16:30'
$ws1.Range("F10").Value = 380
$ws1.Range("G10").Value = 50
$ws1.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=85000'
$ws1.Range("I10").Value = '//i2.hdslb.com/bfs/openplatform/202404/U2EZscfQ1714448575403.jpeg'

$ws1.Range("B11").Value = '2024-07-20'
$ws1.Range("C11").Value = '合肥·W·A首届童年怀旧only'
$ws1.Range("D11").Value = '铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)'
$ws1.Range("E11").Value = '2024.07.20 09:30-07.20 17:00'
$ws1.Range("F11").Value = 122
$ws1.Range("G11").Value = 78
$ws1.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=84794'
$ws1.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202404/Ie0KTNEr1713951888990.png'

$ws1.Range("B12").Value = '2024-07-27'
$ws1.Range("C12").Value = '安徽·MAX特摄only展'
$ws1.Range("D12").Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws1.Range("E12").Value = '2024.07.27 09:30-07.27 18:00'
$ws1.Range("F12").Value = 186
$ws1.Range("G12").Value = 50
$ws1.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=83684'
$ws1.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202404/jv1CiqqW1712029200830.jpeg'

$ws1.Range("B13").Value = '2024-08-03'
$ws1.Range("C13").Value = '合肥·第七届环形宇宙动漫游戏嘉年华'
$ws1.Range("D13").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws1.Range("E13").Value = '2024.08.03 09:30-08.04 17:00'
$ws1.Range("F13").Value = 555
$ws1.Range("G13").Value = 49
$ws1.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=84767'
$ws1.Range("I13").Value = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'

# Remove now-stale trailing rows (shrinks dimension to A1:I13)
for ($i = 0; $i -lt 3; $i++) {
    $ws1.Rows.Item(14).Delete()
}

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("B2").Value = '2024-05-18'
$ws4.Range("C2").Value = '合肥·WA二次元饭局（取消）'
$ws4.Range("D2").Value = '临泉路胜利路交叉路（中环国际大厦对面） 太太满庭芳(胜利路店)'
$ws4.Range("E2").Value = '2024.05.18 14:50-05.18 20:00'
$ws4.Range("F2").Value = 67
$ws4.Range("G2").Value = '不可售'
$ws4.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=83978'
$ws4.Range("I2").Value = '//i1.hdslb.com/bfs/openplatform/202404/wK9Yq9Ta1712657384067.jpeg'

$ws4.Range("B3").Value = '2024-05-18'
$ws4.Range("C3").Value = '合肥·梦时空SPO1动漫展（取消）'
$ws4.Range("D3").Value = '阜阳路16号 银瑞林国际大酒店'
$ws4.Range("E3").Value = '2024.05.18 10:00-05.18 17:00'
$ws4.Range("F3").Value = 127
$ws4.Range("G3").Value = '不可售'
$ws4.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
$ws4.Range("I3").Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'

$ws4.Range("B4").Value = '2024-05-18'
$ws4.Range("C4").Value = '合肥·首届偶活企划——偶像计划-闪耀舞台'
$ws4.Range("D4").Value = '阜阳路16号 银瑞林国际大酒店'
$ws4.Range("E4").Value = '2024.05.18 09:00-05.18 17:00'
$ws4.Range("F4").Value = 60
$ws4.Range("G4").Value = 58
$ws4.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=83891'
$ws4.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202404/lfqv8l9Q1712453982625.jpeg'

$ws4.Range("B5").Value = '2024-06-01'
$ws4.Range("C5").Value = '合肥·运动番only·群青日和'
$ws4.Range("D5").Value = '金寨路287号 合肥明星运动公园'
$ws4.Range("E5").Value = '2024.06.01 09:30-06.01 17:30'
$ws4.Range("F5").Value = 504
$ws4.Range("G5").Value = 70
$ws4.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=83058'
$ws4.Range("I5").Value = '//i2.hdslb.com/bfs/openplatform/202404/Jzeq47lD1714026878824.jpeg'

$ws4.Range("B6").Value = '2024-06-08'
$ws4.Range("C6").Value = '合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~'
$ws4.Range("D6").Value = '锦绣大道3899号 合肥滨湖会展中心'
$ws4.Range("E6").Value = '2024.06.08 09:30-06.09 17:00'
$ws4.Range("F6").Value = 6740
$ws4.Range("G6").Value = 65
$ws4.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=83518'
$ws4.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202403/1Sqp42gM1711691520194.jpeg'

$ws4.Range("B7").Value = '2024-06-09'
$ws4.Range("C7").Value = '合肥·第二届华盟动漫次元嘉年华'
$ws4.Range("D7").Value = '常青街道十五里河村合柴1972院内 合肥当代美术馆'
$ws4.Range("E7").Value = '2024.06.09 10:00-06.10 17:00'
$ws4.Range("F7").Value = 186
$ws4.Range("G7").Value = 58
$ws4.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=84081'
$ws4.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202404/O5LyHE7j1712732240786.jpeg'

$ws4.Range("B8").Value = '2024-06-09'
$ws4.Range("C8").Value = '合肥·第六届环形宇宙动漫游戏嘉年华内场票·赵成晨'
$ws4.Range("D8").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws4.Range("E8").Value = '2024.06.09 09:30-06.09 17:00'
$ws4.Range("F8").Value = 148
$ws4.Range("G8").Value = 238
$ws4.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=84863'
$ws4.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202404/I5S4Ih2M1714031127805.jpeg'

$ws4.Range("B9").Value = '2024-06-22'
$ws4.Range("C9").Value = '合肥·Look Look动漫嘉年华'
$ws4.Range("D9").Value = '新站区东方大道288号 少荃体育中心'
$ws4.Range("E9").Value = '2024.06.22 10:00-06.22 17:30'
$ws4.Range("F9").Value = 1035
$ws4.Range("G9").Value = 58
$ws4.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=82311'
$ws4.Range("I9").Value = '//i2.hdslb.com/bfs/openplatform/202404/RFYwkzvt1713951750482.jpeg'

$ws4.Range("B10").Value = '2024-06-22'
$ws4.Range("C10").Value = '合肥·城市动漫节'
$ws4.Range("D10").Value = '包河经济开发区大连路与园博大道交口骆岗中央公园园博小镇一期S6区1号楼 大机库演艺中心'
$ws4.Range("E10").Value = '2024.06.22 10:00-06.23 16:30'
$ws4.Range("F10").Value = 380
$ws4.Range("G10").Value = 50
$ws4.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=85000'
$ws4.Range("I10").Value = '//i2.hdslb.com/bfs/openplatform/202404/U2EZscfQ1714448575403.jpeg'

$ws4.Range("B11").Value = '2024-07-20'
$ws4.Range("C11").Value = '合肥·W·A首届童年怀旧only'
$ws4.Range("D11").Value = '铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)'
$ws4.Range("E11").Value = '2024.07.20 09:30-07.20 17:00'
$ws4.Range("F11").Value = 122
$ws4.Range("G11").Value = 78
$ws4.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=84794'
$ws4.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202404/Ie0KTNEr1713951888990.png'

$ws4.Range("B12").Value = '2024-07-27'
$ws4.Range("C12").Value = '安徽·MAX特摄only展'
$ws4.Range("D12").Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws4.Range("E12").Value = '2024.07.27 09:30-07.27 18:00'
$ws4.Range("F12").Value = 186
$ws4.Range("G12").Value = 50
$ws4.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=83684'
$ws4.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202404/jv1CiqqW1712029200830.jpeg'

$ws4.Range("B13").Value = '2024-08-03'
$ws4.Range("C13").Value = '合肥·第七届环形宇宙动漫游戏嘉年华'
$ws4.Range("D13").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws4.Range("E13").Value = '2024.08.03 09:30-08.04 17:00'
$ws4.Range("F13").Value = 555
$ws4.Range("G13").Value = 49
$ws4.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=84767'
$ws4.Range("I13").Value = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'

$ws4.Range("B14").Value = '2024-08-03'
$ws4.Range("C14").Value = '合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会'
$ws4.Range("D14").Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range("E14").Value = '2024.08.03 19:30-08.03 21:00'
$ws4.Range("F14").Value = 18
$ws4.Range("G14").Value = 80
$ws4.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=83556'
$ws4.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg'

# Remove now-stale trailing rows (shrinks dimension to A1:I14)
for ($i = 0; $i -lt 3; $i++) {
    $ws4.Rows.Item(15).Delete()
}

Write-Output "edit complete"